$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(51, 28, "PAMR4",    2,  5805.515239477503, 1),
    @(52, 28, "CAND234",  6,  5381.165919282511, 0),
    @(53, 28, "PMP234",  12, 29352.43049327355,  0),
    @(54, 22, "PAMR2234", 5,   79.7457627118644, 0)
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
